# Atualização de bases das ligas, do dia: 19-06-2024 às 21:51
#
# The underlying source data had several rows whose records got re-sorted
# (the "id" sequence in column A stays fixed per row, but the match record
# that belongs with each sequence number changed). For each of the pairs
# below, the two rows simply trade their entire record (every column from
# B "id" through AD "PL_AhUnder") while column A (the running index) stays
# put on its own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(34, 35),
    @(126, 127),
    @(130, 131),
    @(132, 133),
    @(175, 176),
    @(322, 323)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $v1 = $range1.Value2
    $v2 = $range2.Value2

    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

Write-Output "Swapped $($rowPairs.Count) row pairs"
